# Actualización automática hashcode lun mar 15 04:58:28 CET 2021
# Updates the hashcode values (column B) for a set of rows identified by
# their code in column A, matching the shared-string text replacements
# described by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hashcode.csv")

$updates = @{
    "B44"  = "38d8a5eb271cd589ddc82f72caa28981"
    "B89"  = "a55aa003b040b5d58bd2eb051d227a34"
    "B154" = "03cb6b26fa705d888bc739a8fef89740"
    "B161" = "f624c541439ac55856f2147391a085ff"
    "B223" = "88f9769d853035477b76e80a90137177"
    "B226" = "6be3c75043c00b526393ce67c30f59a1"
    "B231" = "f78871005dfa8b1ce50e3cdf572e411a"
    "B248" = "497a413a8d04100d6e3242bc7d878332"
    "B282" = "fcb00ade864bb53decdbaab35c770ebf"
    "B417" = "58fcbf20ee5a33e6680ed4703e3e18ea"
    "B454" = "64a227d433778bf34cf8a8fd975face6"
    "B530" = "a4d2568039fb421909de6d201afbd562"
    "B569" = "34af0eaefed423f6e7d25dbb1f86b0a9"
    "B582" = "47ea736a6a77c2fc6434aeac7ae84d4e"
    "B772" = "78a8bcccab684fb25543376660854a81"
    "B778" = "07fd5ff28a5d01d0bb6287e70e7075a6"
    "B803" = "f4dfc40fb110a1d58bc00ccb8747d2ca"
    "B844" = "f678ee8185fc63e3bb2bc5a4070e68c4"
    "B873" = "1e34f80982028ed80cafa842d95a8b1c"
    "B887" = "6b9a4742a2d9f43aafb4c82b3acc4602"
    "B892" = "d878f735a89572d2273c1e98708e28dd"
    "B923" = "a9d84eeeb25cd8afeac1db1fbfcf6e9c"
    "B947" = "22b87dd39c5a5cdc21c84e5bb5b0fa4b"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
